# schematic.pptx edit:
#   * Methods text fix: "... even though four sequences have the mutation."
#                     -> "... although more than three sequences have the mutation."
#     (shape "TextBox 43" / creationId {78B8061B-22FC-C02B-A7E3-6982D9CFC56E},
#      nested inside the top-level "Group 39" on slide 1)
#   * Footer date placeholder re-cache 1/12/23 -> 1/27/23 on the slide master and
#     every slide layout (PowerPoint stamps the `datetimeFigureOut` field with the
#     save date whenever the deck is touched/re-saved).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix the Methods slide text describing how mutations are counted.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$group = $s.Shapes.Item(1)
$items = $group.GroupItems

$target = $null
for ($i = 1; $i -le $items.Count; $i++) {
    $sh = $items.Item($i)
    if ($sh.Id -eq 44) {
        $target = $sh
        break
    }
}

$tr = $target.TextFrame.TextRange
$fullText = $tr.Text
$oldPhrase = " mutations even though four sequences have the mutation."
$newPhrase = " mutations although more than three sequences have the mutation."
$startIdx = $fullText.IndexOf($oldPhrase)
if ($startIdx -ge 0) {
    $chars = $tr.Characters($startIdx + 1, $oldPhrase.Length)
    $chars.Text = $newPhrase
}

# ---------------------------------------------------------------------------
# 2) Re-cache the "datetimeFigureOut" footer field text on the slide master
#    and on every custom (slide) layout, 1/12/23 -> 1/27/23.
# ---------------------------------------------------------------------------
$oldDate = "1/12/23"
$newDate = "1/27/23"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*" -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
